{"js": "// Convert the field code `{ m:'doc.html'.fromHTMLURI() }` (stored as a real\n// Word field: fldChar begin / instrText* / fldChar end) into plain literal\n// text runs reading `{m:'doc.html'.fromHTMLURI()}`, while keeping the\n// `_GoBack` bookmark exactly where it was (between \"doc.html\" and\n// \"'.fromHTMLURI()\").\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Locate the paragraph that actually holds the field (robust to position).\nfor (const p of paragraphs.items) {\n  p.fields.load(\"items\");\n}\nawait context.sync();\n\nlet target = null;\nfor (const p of paragraphs.items) {\n  if (p.fields.items.length > 0) {\n    target = p;\n    break;\n  }\n}\n\nif (target) {\n  const range = target.getRange();\n\n  // Replace the whole paragraph's contents with plain-text runs that carry\n  // the same literal characters the field code used to contain, keeping the\n  // bookmark in place between \"doc.html\" and \"'.fromHTMLURI()\".\n  const ooxml =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\" ' +\n    'pkg:padding=\"512\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n    '<w:p>' +\n    '<w:r><w:t>{</w:t></w:r>' +\n    '<w:r><w:t>m</w:t></w:r>' +\n    '<w:r><w:t>:</w:t></w:r>' +\n    \"<w:r><w:t>'</w:t></w:r>\" +\n    '<w:r><w:t>doc.html</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n    '<w:bookmarkEnd w:id=\"0\"/>' +\n    \"<w:r><w:t>'.fromHTMLURI()</w:t></w:r>\" +\n    '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>' +\n    '</w:p>' +\n    '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>';\n\n  range.insertOoxml(ooxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Convert the field code `{ m:'doc.html'.fromHTMLURI() }` (stored as a real\n# Word field: fldChar begin / instrText* / fldChar end) into plain literal\n# text runs reading `{m:'doc.html'.fromHTMLURI()}`, while keeping the\n# `_GoBack` bookmark exactly where it was (between \"doc.html\" and\n# \"'.fromHTMLURI()\").\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that actually holds the field (robust to position).\n$targetPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Fields.Count -gt 0) {\n        $targetPara = $p\n        break\n    }\n}\n\nif ($targetPara -ne $null) {\n    # Range covering the paragraph's content but excluding its trailing\n    # paragraph mark, so the replacement stays inside this paragraph.\n    $r = $d.Range($targetPara.Range.Start, $targetPara.Range.End - 1)\n\n    $apos = \"'\"\n    $ooxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" ' +\n        'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\" ' +\n        'pkg:padding=\"512\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' +\n        '<w:p>' +\n        '<w:r><w:t>{</w:t></w:r>' +\n        '<w:r><w:t>m</w:t></w:r>' +\n        '<w:r><w:t>:</w:t></w:r>' +\n        '<w:r><w:t>' + $apos + '</w:t></w:r>' +\n        '<w:r><w:t>doc.html</w:t></w:r>' +\n        '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n        '<w:bookmarkEnd w:id=\"0\"/>' +\n        '<w:r><w:t>' + $apos + '.fromHTMLURI()</w:t></w:r>' +\n        '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>' +\n        '</w:p>' +\n        '</w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData>' +\n        '</pkg:part>' +\n        '</pkg:package>'\n\n    [void]$r.InsertXML($ooxml)\n}\n"}
